{"js": "// Applies the wording/copy-edit pass described by the diff:\n//  - \"a customer's table\"            -> \"a customers' table\"\n//  - \"information of each\"           -> \"information on each\"\n//  - \"our Python apps) ... for them\" -> \"our Python app) ... from them\"\n//  - \"(your applications)\"           -> \"(applications)\"\n//  - \"do not have to use\"            -> \"do not need to use\"\n//  - \"is related then expect\"        -> \"is related, then expect\"\n//  - \"the applications grow and\"     -> \"the applications' growth and\"\n// and relocates the (Word-managed) \"_GoBack\" last-edit bookmark from the\n// intro paragraph to right after the newly-typed word \"growth\" - exactly\n// where Word itself would drop it after that keystroke.\n\nconst doc = context.document;\nconst body = doc.body;\n\nasync function replaceOnce(searchText, newText, options) {\n  const results = body.search(searchText, Object.assign({ matchCase: true }, options));\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\n      \"Expected exactly 1 match for \" + JSON.stringify(searchText) +\n      \" but found \" + results.items.length\n    );\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1. \"We have three tables ... a customer's table\" -> \"... a customers' table\"\nawait replaceOnce(\"customer\\u2019s table\", \"customers\\u2019 table\");\n\n// 2. \"the items' table defines information of each\" -> \"... information on each\"\nawait replaceOnce(\"information of each\", \"information on each\");\n\n// 3. \"Applications (such as our Python apps) ...\" -> \"... Python app) ...\"\nawait replaceOnce(\"such as our Python apps)\", \"such as our Python app)\");\n\n// 4. \"... ask it to retrieve data for them\" -> \"... data from them\"\nawait replaceOnce(\"data for them\", \"data from them\");\n\n// 5. \"... serve more clients (your applications)\" -> \"... (applications)\"\nawait replaceOnce(\"your applications)\", \"applications)\");\n\n// 6. \"... then we do not have to use a relational system\" -> \"... need to use ...\"\nawait replaceOnce(\"do not have to use\", \"do not need to use\");\n\n// 7. \"If our data is related then expect that as\" -> \"..., then expect that as\"\nawait replaceOnce(\"is related then expect\", \"is related, then expect\");\n\n// 8. \"the applications grow and more data\" -> \"the applications' growth and more data\"\nawait replaceOnce(\"the applications grow and\", \"the applications\\u2019 growth and\");\n\n// 9. Word keeps a single \"_GoBack\" bookmark marking the last edit location.\n//    It used to sit in the intro paragraph; after this edit session it\n//    belongs right after the word we just finished typing (\"growth\").\ndoc.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst growthResults = body.search(\"growth\", { matchCase: true });\ngrowthResults.load(\"text\");\nawait context.sync();\nif (growthResults.items.length !== 1) {\n  throw new Error(\"Expected exactly 1 match for \\\"growth\\\" but found \" + growthResults.items.length);\n}\ngrowthResults.items[0].getRange(Word.RangeLocation.end).insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Applies the wording/copy-edit pass described by the diff:\n#  - \"a customer's table\"            -> \"a customers' table\"\n#  - \"information of each\"           -> \"information on each\"\n#  - \"our Python apps) ... for them\" -> \"our Python app) ... from them\"\n#  - \"(your applications)\"           -> \"(applications)\"\n#  - \"do not have to use\"            -> \"do not need to use\"\n#  - \"is related then expect\"        -> \"is related, then expect\"\n#  - \"the applications grow and\"     -> \"the applications' growth and\"\n# and relocates the (Word-managed) \"_GoBack\" last-edit bookmark from the\n# intro paragraph to right after the newly-typed word \"growth\" - exactly\n# where Word itself would drop it after that keystroke.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Once {\n    param(\n        [string]$SearchText,\n        [string]$ReplaceText\n    )\n\n    # Make sure we only ever touch a single, unambiguous occurrence.\n    $fullText = $d.Content.Text\n    $escaped = [regex]::Escape($SearchText)\n    $matches = [regex]::Matches($fullText, $escaped)\n    if ($matches.Count -ne 1) {\n        throw \"Expected exactly 1 match for '$SearchText' but found $($matches.Count)\"\n    }\n\n    $rng = $d.Content\n    $found = $rng.Find.Execute($SearchText, $true, $false, $false, $false, $false, $true, 1, $false, $ReplaceText, 2)\n    if (-not $found) {\n        throw \"Find.Execute failed to replace '$SearchText'\"\n    }\n}\n\n# 1. \"We have three tables ... a customer's table\" -> \"... a customers' table\"\nReplace-Once \"customer\u2019s table\" \"customers\u2019 table\"\n\n# 2. \"the items' table defines information of each\" -> \"... information on each\"\nReplace-Once \"information of each\" \"information on each\"\n\n# 3. \"Applications (such as our Python apps) ...\" -> \"... Python app) ...\"\nReplace-Once \"such as our Python apps)\" \"such as our Python app)\"\n\n# 4. \"... ask it to retrieve data for them\" -> \"... data from them\"\nReplace-Once \"data for them\" \"data from them\"\n\n# 5. \"... serve more clients (your applications)\" -> \"... (applications)\"\nReplace-Once \"your applications)\" \"applications)\"\n\n# 6. \"... then we do not have to use a relational system\" -> \"... need to use ...\"\nReplace-Once \"do not have to use\" \"do not need to use\"\n\n# 7. \"If our data is related then expect that as\" -> \"..., then expect that as\"\nReplace-Once \"is related then expect\" \"is related, then expect\"\n\n# 8. \"the applications grow and more data\" -> \"the applications' growth and more data\"\nReplace-Once \"the applications grow and\" \"the applications\u2019 growth and\"\n\n# 9. Word keeps a single \"_GoBack\" bookmark marking the last edit location.\n#    It used to sit in the intro paragraph; after this edit session it\n#    belongs right after the word we just finished typing (\"growth\").\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n$growthRange = $d.Content\n$found = $growthRange.Find.Execute(\"growth\")\nif (-not $found) {\n    throw \"Could not find 'growth' to place the _GoBack bookmark\"\n}\n$growthRange.Collapse(0)  # wdCollapseEnd\n$d.Bookmarks.Add(\"_GoBack\", $growthRange) | Out-Null\n"}
